$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "30.219.76"
$ws.Range("E2").Value = "  -0.47%  "
$ws.Range("D3").Value = "1.861.69"
$ws.Range("E3").Value = "  -0.51%  "
$ws.Range("D4").Value = "'1.001"
$ws.Range("D5").Value = "'236.51"
$ws.Range("E5").Value = "  +0.41%  "
$ws.Range("D6").Value = "'1.000"
$ws.Range("E6").Value = "  -0.11%  "
$ws.Range("D7").Value = "'0.4683"
$ws.Range("E7").Value = "  +0.38%  "
$ws.Range("D8").Value = "'0.2898"
$ws.Range("E8").Value = "  +2.13%  "
$ws.Range("D9").Value = "'0.06535"
$ws.Range("E9").Value = "  -0.25%  "
$ws.Range("D10").Value = "'21.64"
$ws.Range("E10").Value = "  +2.63%  "
$ws.Range("D11").Value = "'0.07932"
$ws.Range("E11").Value = "  -0.03%  "
$ws.Range("D12").Value = "'98.11"
$ws.Range("E12").Value = "  +0.87%  "
$ws.Range("D13").Value = "1.868.85"
$ws.Range("E13").Value = "  -0.42%  "
$ws.Range("D14").Value = "'5.166"
$ws.Range("E14").Value = "  +0.23%  "
$ws.Range("E15").Value = "  +0.86%  "
$ws.Range("D16").Value = "'267.56"
$ws.Range("E16").Value = "  -5.05%  "
$ws.Range("D17").Value = "30.212.78"
$ws.Range("E17").Value = "  -0.52%  "
$ws.Range("E18").Value = "  +8.67%  "
$ws.Range("D19").Value = "'1.000"
$ws.Range("E19").Value = "  -0.03%  "
$ws.Range("D20").Value = "'0.000007383"
$ws.Range("E20").Value = "  +1.28%  "
$ws.Range("D21").Value = "2.114.41"
$ws.Range("E21").Value = "  -0.29%  "
$ws.Range("D22").Value = "'5.317"
$ws.Range("E22").Value = "  -4.12%  "
$ws.Range("D23").Value = "'1.000"
$ws.Range("E23").Value = "  -0.10%  "
$ws.Range("D24").Value = "'6.179"
$ws.Range("E24").Value = "  -0.30%  "
$ws.Range("D25").Value = "'166.58"
$ws.Range("E25").Value = "  +1.30%  "
$ws.Range("D26").Value = "'9.213"
$ws.Range("E26").Value = "  -0.83%  "
$ws.Range("E27").Value = "  -1.07%  "
$ws.Range("E28").Value = "  +1.19%  "
$ws.Range("D29").Value = "'1.390"
$ws.Range("E29").Value = "  +2.64%  "
$ws.Range("D30").Value = "'0.09845"
$ws.Range("E30").Value = "  +1.61%  "
$ws.Range("D31").Value = "'4.369"
$ws.Range("E31").Value = "  -1.59%  "
$ws.Range("D32").Value = "'1.470"
$ws.Range("E32").Value = "  -0.40%  "
$ws.Range("D33").Value = "'4.049"
$ws.Range("E33").Value = "  -1.64%  "
$ws.Range("D34").Value = "'0.04704"
$ws.Range("E34").Value = "  -0.13%  "
$ws.Range("D35").Value = "'1.130"
$ws.Range("E35").Value = "  +0.93%  "
$ws.Range("D36").Value = "'0.7030"
$ws.Range("E36").Value = "  -0.31%  "
$ws.Range("D37").Value = "'2.706"
$ws.Range("E37").Value = "  -0.39%  "
$ws.Range("D38").Value = "'0.01871"
$ws.Range("E38").Value = "  +0.63%  "
$ws.Range("D39").Value = "'2.614"
$ws.Range("E39").Value = "  +2.74%  "
$ws.Range("D40").Value = "'6.292"
$ws.Range("E40").Value = "  -0.75%  "
$ws.Range("D41").Value = "'74.24"
$ws.Range("E41").Value = "  +0.66%  "
$ws.Range("D42").Value = "'1.935"
$ws.Range("E42").Value = "  -0.60%  "
$ws.Range("D43").Value = "'0.8459"
$ws.Range("E43").Value = "  -0.29%  "
$ws.Range("B44").Value = "TheSandbox"
$ws.Range("C44").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D44").Value = "'0.4161"
$ws.Range("E44").Value = "  -0.81%  "
$ws.Range("B45").Value = "PaxDollar"
$ws.Range("C45").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D45").Value = "'0.9994"
$ws.Range("E45").Value = "  -0.15%  "
$ws.Range("D46").Value = "'102.91"
$ws.Range("E46").Value = "  -0.94%  "
$ws.Range("D47").Value = "'956.01"
$ws.Range("E47").Value = "  +2.03%  "
$ws.Range("D48").Value = "'7.148"
$ws.Range("E48").Value = "  -0.90%  "
$ws.Range("D49").Value = "'9.233"
$ws.Range("E49").Value = "  -0.73%  "
$ws.Range("D50").Value = "'34.11"
$ws.Range("E50").Value = "  -0.19%  "
$ws.Range("B51").Value = "Decentraland"
$ws.Range("C51").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D51").Value = "'0.3898"
$ws.Range("E51").Value = "  +2.24%  "
